$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '47.380.18'
$ws.Range('E2').Value = '  +4.29%  '
$ws.Range('D3').Value = '2.487.00'
$ws.Range('E3').Value = '  +2.49%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''322.90'
$ws.Range('E5').Value = '  +1.34%  '
$ws.Range('D6').Value = '''106.59'
$ws.Range('E6').Value = '  +3.43%  '
$ws.Range('E7').Value = '  +1.62%  '
$ws.Range('D8').Value = '''0.999'
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  +2.44%  '
$ws.Range('D10').Value = '''38.24'
$ws.Range('E10').Value = '  +7.24%  '
$ws.Range('E11').Value = '  +1.05%  '
$ws.Range('E12').Value = '  +1.23%  '
$ws.Range('E13').Value = '  +0.87%  '
$ws.Range('E14').Value = '  +1.22%  '
$ws.Range('D15').Value = '2.876.80'
$ws.Range('E15').Value = '  +2.50%  '
$ws.Range('D16').Value = '2.487.21'
$ws.Range('E16').Value = '  +2.87%  '
$ws.Range('D17').Value = '''0.845'
$ws.Range('E17').Value = '  +0.26%  '
$ws.Range('D18').Value = '47.268.75'
$ws.Range('E18').Value = '  +4.26%  '
$ws.Range('D19').Value = '''12.80'
$ws.Range('E19').Value = '  +4.65%  '
$ws.Range('E20').Value = '  +3.34%  '
$ws.Range('D21').Value = '0.0₃0938'
$ws.Range('E21').Value = '  +1.45%  '
$ws.Range('D22').Value = '''70.65'
$ws.Range('E22').Value = '  +2.39%  '
$ws.Range('D23').Value = '''2.43'
$ws.Range('E23').Value = '  +6.72%  '
$ws.Range('D24').Value = '''251.56'
$ws.Range('E24').Value = '  +2.81%  '
$ws.Range('E25').Value = '  +3.25%  '
$ws.Range('D26').Value = '''26.13'
$ws.Range('E26').Value = '  +1.45%  '
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').Value = '''10.02'
$ws.Range('E28').Value = '  +4.29%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '''2.24'
$ws.Range('E29').Value = '  -2.04%  '
$ws.Range('D30').Value = '''35.22'
$ws.Range('E30').Value = '  +6.91%  '
$ws.Range('E31').Value = '  +8.39%  '
$ws.Range('D32').Value = '''49.48'
$ws.Range('E32').Value = '  +0.10%  '
$ws.Range('D33').Value = '''19.68'
$ws.Range('E33').Value = '  -3.07%  '
$ws.Range('E34').Value = '  +3.39%  '
$ws.Range('D35').Value = '''0.0788'
$ws.Range('E35').Value = '  +2.87%  '
$ws.Range('E36').Value = '  +0.18%  '
$ws.Range('E37').Value = '  +5.45%  '
$ws.Range('D38').Value = '''4.63'
$ws.Range('E38').Value = '  +3.49%  '
$ws.Range('E39').Value = '  +3.56%  '
$ws.Range('E40').Value = '  +1.94%  '
$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').Value = '''122.11'
$ws.Range('E41').Value = '  -2.80%  '
$ws.Range('B42').Value = 'WEMIXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D42').Value = '''2.24'
$ws.Range('E42').Value = '  +1.86%  '
$ws.Range('D43').Value = '''21.08'
$ws.Range('E43').Value = '  +2.21%  '
$ws.Range('E44').Value = '  +2.40%  '
$ws.Range('D45').Value = '1.962.85'
$ws.Range('E45').Value = '  +1.22%  '
$ws.Range('E46').Value = '  +2.01%  '
$ws.Range('E47').Value = '  -0.43%  '
$ws.Range('D48').Value = '''1.80'
$ws.Range('E48').Value = '  +0.61%  '
$ws.Range('E49').Value = '  -0.09%  '
$ws.Range('E50').Value = '  +9.74%  '
$ws.Range('D51').Value = '''79.65'
$ws.Range('E51').Value = '  +3.47%  '
